$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextCell "D2" "26.299.14"
Set-TextCell "E2" "  +1.23%  "
Set-TextCell "D3" "1.620.98"
Set-TextCell "E3" "  +1.71%  "
Set-TextCell "E4" "  +0.05%  "
Set-TextCell "D5" "212.20"
Set-TextCell "E5" "  +0.78%  "
Set-TextCell "E6" "  +0.03%  "
Set-TextCell "E7" "  +0.95%  "
Set-TextCell "D8" "0.248"
Set-TextCell "E8" "  +1.13%  "
Set-TextCell "E9" "  +0.77%  "
Set-TextCell "D10" "18.86"
Set-TextCell "E10" "  +5.13%  "
Set-TextCell "D11" "0.0815"
Set-TextCell "E11" "  +0.77%  "
Set-TextCell "D12" "1.846.84"
Set-TextCell "E12" "  +1.72%  "
Set-TextCell "D13" "1.622.54"
Set-TextCell "E13" "  +1.88%  "
Set-TextCell "E14" "  +0.10%  "
Set-TextCell "E15" "  +1.19%  "
Set-TextCell "D16" "26.306.72"
Set-TextCell "E16" "  +1.24%  "
Set-TextCell "D17" "62.42"
Set-TextCell "E17" "  +3.90%  "
Set-TextCell "E18" "  +0.95%  "
Set-TextCell "E19" "  -0.03%  "
Set-TextCell "D20" "201.79"
Set-TextCell "E20" "  +0.99%  "
Set-TextCell "D21" "4.28"
Set-TextCell "E21" "  +1.30%  "
Set-TextCell "E22" "  +1.10%  "
Set-TextCell "D23" "6.04"
Set-TextCell "E23" "  +0.83%  "
Set-TextCell "E24" "  +4.82%  "
Set-TextCell "D25" "144.03"
Set-TextCell "E25" "  +1.03%  "
Set-TextCell "E26" "  +0.04%  "
Set-TextCell "E27" "  -0.61%  "
Set-TextCell "E28" "  +0.73%  "
Set-TextCell "D29" "6.55"
Set-TextCell "E29" "  +1.76%  "
Set-TextCell "D30" "0.0525"
Set-TextCell "E30" "  +10.64%  "
Set-TextCell "E31" "  +0.82%  "
Set-TextCell "E32" "  +1.72%  "
Set-TextCell "D33" "2.94"
Set-TextCell "E33" "  +0.03%  "
Set-TextCell "E34" "  +1.39%  "
Set-TextCell "E35" "  +2.31%  "
Set-TextCell "D36" "1.179.60"
Set-TextCell "E36" "  +4.92%  "
Set-TextCell "D37" "0.0163"
Set-TextCell "E37" "  +0.53%  "
Set-TextCell "E38" "  +3.17%  "
Set-TextCell "E39" "  +0.08%  "
Set-TextCell "D41" "0.495"
Set-TextCell "E41" "  +1.05%  "
Set-TextCell "D42" "0.787"
Set-TextCell "E42" "  +0.59%  "
Set-TextCell "D43" "5.33"
Set-TextCell "E43" "  +4.68%  "
Set-TextCell "D44" "1.758.22"
Set-TextCell "E44" "  +1.83%  "
Set-TextCell "D45" "93.30"
Set-TextCell "E45" "  +0.80%  "
Set-TextCell "E46" "  +14.27%  "
Set-TextCell "E47" "  +2.05%  "
Set-TextCell "D48" "53.88"
Set-TextCell "E48" "  +1.06%  "
Set-TextCell "E49" "  +1.11%  "
Set-TextCell "E50" "  +0.14%  "
Set-TextCell "E51" "  -0.10%  "
